$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in columns P, Q, R for rows 5-7 ---
$ws.Range("P5").Value = 23.111083656771282
$ws.Range("Q5").Value = 24.08077930418019
$ws.Range("R5").Value = 19.336931533747723

$ws.Range("P6").Value = 14.322631450320875
$ws.Range("Q6").Value = 13.073459110725862
$ws.Range("R6").Value = 10.464141365743002

$ws.Range("P7").Value = 23.612622725489956

# --- Add new column S (year 2022) by copying formats from column R ---
# Row 3: empty bottom-border cell
$ws.Range("R3").Copy() | Out-Null
$ws.Range("S3").PasteSpecial(-4122) | Out-Null

# Row 4: year header
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null
$ws.Range("S4").Value = 2022

# Row 5
$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null
$ws.Range("S5").Value = 13.600365850576139

# Row 6
$ws.Range("R6").Copy() | Out-Null
$ws.Range("S6").PasteSpecial(-4122) | Out-Null
$ws.Range("S6").Value = 9.2742414863791556

# Row 7
$ws.Range("R7").Copy() | Out-Null
$ws.Range("S7").PasteSpecial(-4122) | Out-Null
$ws.Range("S7").Value = 17.303523954725925

# Row 8
$ws.Range("R8").Copy() | Out-Null
$ws.Range("S8").PasteSpecial(-4122) | Out-Null
$ws.Range("S8").Value = 205.5

$excel.CutCopyMode = 0

# --- Update the saved selection to match the author's final cursor position ---
$ws.Range("Q15").Select() | Out-Null
